$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83, shifting rows 83:129 down to 84:130
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the new weekly record
$ws.Cells.Item(83, 1).Value = 11
$ws.Cells.Item(83, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(83, 3).Value = "Bíobío"
$ws.Cells.Item(83, 4).Value = 45141
$ws.Cells.Item(83, 5).Value = 8
$ws.Cells.Item(83, 6).Value = 100112037
$ws.Cells.Item(83, 7).Value = "Cebollín"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 100
$ws.Cells.Item(83, 11).Value = 4500
$ws.Cells.Item(83, 12).Value = 5000
$ws.Cells.Item(83, 13).Value = 4750
$ws.Cells.Item(83, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(83, 15).Value = "Región Metropolitana"
$ws.Cells.Item(83, 16).Value = 132
$ws.Cells.Item(83, 17).Value = 36
$ws.Cells.Item(83, 18).Value = "Hortaliza"
